# Fill in the first blank table row (after "Chandralekha Gude") with a new
# name ("Hunny Keshwani", with spell-check markers around each word) and
# the matching skills ("Java, HTML, CSS"), preserving each paragraph's
# existing identity (w14:paraId/w:rsid*) and the fr-FR language formatting
# already used throughout the table.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the first completely empty row beneath the populated entries
# (defensively, in case the row index ever shifts) instead of assuming a
# fixed row number.
$targetRow = 6
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $n = $t.Cell($i, 1).Range.Text.Trim([char]13, [char]7)
    $s = $t.Cell($i, 2).Range.Text.Trim([char]13, [char]7)
    if ($n -eq "" -and $s -eq "") {
        $targetRow = $i
        break
    }
}

# Name cell -> "Hunny Keshwani"
$nameCell = $t.Cell($targetRow, 1)
$nameXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p w14:paraId="170DB5E5" w14:textId="77777777" w:rsidR="005C2426" w:rsidRPr="00BE429F" w:rsidRDefault="005C2426">' +
  '<w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>Hunny</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>Keshwani</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'
$nameCell.Range.InsertXML($nameXml) | Out-Null

# Skills cell (row 6, column 2) -> "Java, HTML, CSS"
$skillsCell = $t.Cell(6, 2)
$skillsXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p w14:paraId="0407DDB9" w14:textId="77777777" w:rsidR="005C2426" w:rsidRPr="00BE429F" w:rsidRDefault="005C2426">' +
  '<w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>Java, HTML, CSS</w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'
$skillsCell.Range.InsertXML($skillsXml) | Out-Null
